$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.127.23'
$ws.Range("E2").Value = '  +6.20%  '
$ws.Range("D3").Value = '1.718.10'
$ws.Range("E3").Value = '  +4.00%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = "'333.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.03%  '
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = "'0.3689"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("D8").Value = "'49.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.86%  '
$ws.Range("D9").Value = "'0.3348"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.06%  '
$ws.Range("D10").Value = "'1.190"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.48%  '
$ws.Range("D11").Value = "'0.07483"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.51%  '
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = "'6.331"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.53%  '
$ws.Range("D14").Value = "'20.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.53%  '
$ws.Range("D15").Value = "'6.958"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.57%  '
$ws.Range("D16").Value = '1.716.98'
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("D17").Value = "'0.00001079"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.98%  '
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = "'82.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.35%  '
$ws.Range("D20").Value = "'0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").Value = "'16.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.93%  '
$ws.Range("D22").Value = "'6.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.36%  '
$ws.Range("D23").Value = "'13.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.21%  '
$ws.Range("D24").Value = '26.065.16'
$ws.Range("E24").Value = '  +6.09%  '
$ws.Range("D25").Value = "'2.469"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").Value = "'2.465"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.14%  '
$ws.Range("D27").Value = "'151.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("D28").Value = "'1.372"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +14.22%  '
$ws.Range("D29").Value = "'19.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.01%  '
$ws.Range("D30").Value = '1.911.35'
$ws.Range("E30").Value = '  +4.14%  '
$ws.Range("D31").Value = "'129.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.47%  '
$ws.Range("D32").Value = "'4.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.28%  '
$ws.Range("D33").Value = "'5.967"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.24%  '
$ws.Range("D34").Value = "'0.08554"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").Value = "'1.712"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("D36").Value = "'12.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.10%  '
$ws.Range("D37").Value = "'5.364"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.37%  '
$ws.Range("D38").Value = "'0.06244"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.16%  '
$ws.Range("D39").Value = "'0.02298"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.32%  '
$ws.Range("D40").Value = "'8.599"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.47%  '
$ws.Range("D41").Value = "'0.2144"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("D43").Value = "'14.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.41%  '
$ws.Range("D44").Value = "'0.6180"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.79%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = "'3.835"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = "'0.5913"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.51%  '
$ws.Range("D48").Value = "'128.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.89%  '
$ws.Range("D49").Value = "'2.023"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").Value = "'0.07262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.69%  '
$ws.Range("D51").Value = "'77.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.05%  '
